$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 for "total_orders" so the remaining
# metric rows shift down by one.
$ws.Rows.Item(2).Insert()

# Row 2: new "total_orders" metric
$ws.Range("A2").Value = "total_orders"
$ws.Range("B2").Value = 2.0

# Updated metric values (rows shifted down by the inserted row)
$ws.Range("B3").Value = 3414.62
$ws.Range("B4").Value = -87.52
$ws.Range("B5").Value = -475.3
$ws.Range("B6").Value = 39.02
$ws.Range("B7").Value = 72.26
$ws.Range("B8").Value = 194.88

# New row for "tips" (previously the last row, now pushed to row 9)
$ws.Range("A9").Value = "tips"
$ws.Range("B9").Value = 0.0
